$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data table header / rows: Manager ID based rows -> Report Date based rows ---
$ws.Range("B7").Value  = "Reports in the system"
$ws.Range("A8").Value  = "Report Date"
$ws.Range("A9").Value  = "1/3/2016-7/3/2016"
$ws.Range("A10").Value = "1/3/2016-21/3/2016"
$ws.Range("A11").Value = "1/4/2016 - 7/4/2016"
$ws.Range("B11").Value = "No"

$a12 = @"
Yes- the report exist in the DB
No- The report not Exist in the DB
"@
$ws.Range("A12").Value = $a12

# --- Use case 1 (SuccessfulPresentation) ---
$b16 = @"
open "Reports Page"
Enter The Dates: "1/3/2016-7/3/2016"
press "view weekly Report" 
"@
$ws.Range("B16").Value = $b16

# --- Old "IncorrectID" use case row removed -> row cleared out ---
$ws.Range("A17").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""

# --- MissingManagerID -> MissingDate ---
$ws.Range("A19").Value = "MissingDate"
$b19 = @"
open "Reports Page"
Enter The Dates ""
Press "view weekly Report"

"@
$ws.Range("B19").Value = $b19
$ws.Range("C19").Value = "System throws ""Enter Date"" message"

# --- ManagerNotFound -> ReportNotFound ---
$ws.Range("A20").Value = "ReportNotFound"
$b20 = @"
open "Reports Page"
Enter The Dates "1/4/2016-7/4/2016"
press "view weekly Report"
"@
$ws.Range("B20").Value = $b20
$ws.Range("C20").Value = "System throws ""Report Not found"" message"
$ws.Range("D20").Value = "Report Not Created Yet"

# --- Row heights ---
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 45

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 17.71

# --- Sheet view: scroll position & selection ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("D20").Select()
